$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename header columns:
#   I1: "Then_Question" -> "Then_Goto"
#   J1: "Else_Question" -> "Else_Goto"
$ws.Cells.Item(1, 9).Value = "Then_Goto"
$ws.Cells.Item(1, 10).Value = "Else_Goto"

# Update the active cell selection to I1 (matches the saved view state)
[void]$ws.Range("I1").Select()
